$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ENVIRONMENT"
$ws.Range("A2").Value = "DEV"
$ws.Range("A3").Value = "STAGING"
$ws.Range("A4").Value = "FT"
$ws.Range("A5").Value = "SIT"
$ws.Range("A6").Value = "PERF"
$ws.Range("A7").Value = "PROD"

$ws.Range("B1").Value = "USERNAME"
$ws.Range("C1").Value = "PASSWORD"

$ws.Range("B2").Value = "dev_johnDoe"
$ws.Range("B4").Value = "ft_ryanBlake"
$ws.Range("C2").Value = "password"
$ws.Range("B5").Value = "sit_minaFlyn"
$ws.Range("B6").Value = "perf_younesEr"
$ws.Range("B7").Value = "prof_victoriaU"
$ws.Range("B3").Value = "stage_samMorris"

$ws.Range("C3").Value = "password"
$ws.Range("C4").Value = "password"
$ws.Range("C5").Value = "password"
$ws.Range("C6").Value = "password"
$ws.Range("C7").Value = "password"

# Update the selected cell to match the saved view state
[void]$ws.Range("B12").Select()
